# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# New header cells AD1:AF1 get the same style as the existing header row
# (bold font, thin border, centered/top aligned), and every data row
# (2-42) gets a constant record of 88 wins, 75 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the existing header cells (style index 1:
# bold, thin border, centered/top-aligned) by copying format from the
# last existing header cell onto the new ones.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the record for every player row.
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 75
    $ws.Cells.Item($r, 32).Value = 0
}
